$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = [double]"39.7741120585799"
$ws.Range("D2").Value = [double]"-4446.38897826738"
$ws.Range("E2").Value = [double]"8976.28504729307"
$ws.Range("G2").Value = [double]"1"
$ws.Range("H2").Value = [double]"9186.00759785883"
$ws.Range("I2").Value = [double]"212197.971105761"
$ws.Range("J2").Value = [double]"1082.22588794142"

# Row 3
$ws.Range("F3").Value = [double]"102.961481323311"
$ws.Range("G3").Value = [double]"0.0000000000000000000000438731085200036"

# Row 4
$ws.Range("F4").Value = [double]"103.312436643288"
$ws.Range("G4").Value = [double]"0.0000000000000000000000368120011498852"

# Row 5
$ws.Range("F5").Value = [double]"211.808694838493"
$ws.Range("G5").Value = [double]"0.000000000000000000000000000000000000000000000101467347220903"

# Row 6
$ws.Range("F6").Value = [double]"275.685907498342"
$ws.Range("G6").Value = [double]"0.00000000000000000000000000000000000000000000000000000000000136636213268829"

# Row 7
$ws.Range("F7").Value = [double]"685.960112921679"
$ws.Range("G7").Value = [double]"0.0000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000111084656440435"

# Row 8
$ws.Range("F8").Value = [double]"712.820960688447"
$ws.Range("G8").Value = [double]"0.0000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000163265752374575"

# Row 9
$ws.Range("F9").Value = [double]"718.762431237912"
$ws.Range("G9").Value = [double]"0.000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000836991710901724"

# Row 10
$ws.Range("F10").Value = [double]"770.729641775799"
$ws.Range("G10").Value = [double]"0.00000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000434695124038125"

# Row 11
$ws.Range("F11").Value = [double]"1251.69234825747"
$ws.Range("G11").Value = [double]"1.57928334597889e-272"

# Row 12
$ws.Range("F12").Value = [double]"1492.82395076031"
$ws.Range("G12").Value = [double]"0"
